# Auto-generated Excel COM-interop script applying the Jenova_Profits diff
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H18").Value = 1265
$ws_ALC.Range("I18").Value = 400
$ws_ALC.Range("J18").Value = 2995
$ws_ALC.Range("K18").Value = 400
$ws_ALC.Range("L18").Value = 2995
$ws_ALC.Range("M18").Value = -116
$ws_ALC.Range("N18").Value = -3563
$ws_ALC.Range("H112").Value = 2399.6
$ws_ALC.Range("J112").Value = 2399.6
$ws_ALC.Range("L112").Value = 7198.799999999999
$ws_ALC.Range("N112").Value = -9414.799999999999
$ws_ALC.Range("H125").Value = 5249.2104
$ws_ALC.Range("I125").Value = 5741.375
$ws_ALC.Range("J125").Value = 4891.273
$ws_ALC.Range("K125").Value = 51672.375
$ws_ALC.Range("L125").Value = 44021.457
$ws_ALC.Range("M125").Value = -49212.375
$ws_ALC.Range("N125").Value = -48941.457
$ws_ALC.Range("H134").Value = 58634.816
$ws_ALC.Range("J134").Value = 58634.816
$ws_ALC.Range("L134").Value = 58634.816
$ws_ALC.Range("N134").Value = -68774.81599999999
$ws_ALC.Range("H135").Value = 3119.8215
$ws_ALC.Range("I135").Value = 2894.24
$ws_ALC.Range("K135").Value = 26048.16
$ws_ALC.Range("M135").Value = -23513.16
$ws_ALC.Range("H137").Value = 9024.597
$ws_ALC.Range("I137").Value = 1491.25
$ws_ALC.Range("K137").Value = 4473.75
$ws_ALC.Range("M137").Value = -1923.75
$ws_ALC.Range("H138").Value = 4622.636
$ws_ALC.Range("I138").Value = 1479.8148
$ws_ALC.Range("J138").Value = 6319.76
$ws_ALC.Range("K138").Value = 4439.4444
$ws_ALC.Range("L138").Value = 18959.28
$ws_ALC.Range("M138").Value = 700.5555999999997
$ws_ALC.Range("N138").Value = -29239.28

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 2338.4707
$ws_ARM.Range("I32").Value = 2430.8367
$ws_ARM.Range("J32").Value = 75.5
$ws_ARM.Range("K32").Value = 2430.8367
$ws_ARM.Range("L32").Value = 75.5
$ws_ARM.Range("M32").Value = -2143.8367
$ws_ARM.Range("N32").Value = -649.5
$ws_ARM.Range("H74").Value = 201380.75
$ws_ARM.Range("I74").Value = 218857.94
$ws_ARM.Range("J74").Value = 40590.6
$ws_ARM.Range("K74").Value = 218857.94
$ws_ARM.Range("L74").Value = 40590.6
$ws_ARM.Range("M74").Value = -217983.94
$ws_ARM.Range("N74").Value = -42338.6
$ws_ARM.Range("H77").Value = 201380.75
$ws_ARM.Range("I77").Value = 218857.94
$ws_ARM.Range("J77").Value = 40590.6
$ws_ARM.Range("K77").Value = 1094289.7
$ws_ARM.Range("L77").Value = 202953
$ws_ARM.Range("M77").Value = -1089921.7
$ws_ARM.Range("N77").Value = -211689
$ws_ARM.Range("H132").Value = 261959.62
$ws_ARM.Range("I132").Value = 432640.25
$ws_ARM.Range("K132").Value = 1297920.75
$ws_ARM.Range("M132").Value = -1295390.75

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H81").Value = 58080
$ws_BSM.Range("I81").Value = 30000
$ws_BSM.Range("J81").Value = 62760
$ws_BSM.Range("K81").Value = 30000
$ws_BSM.Range("L81").Value = 62760
$ws_BSM.Range("M81").Value = -28939
$ws_BSM.Range("N81").Value = -64882
$ws_BSM.Range("H84").Value = 58080
$ws_BSM.Range("I84").Value = 30000
$ws_BSM.Range("J84").Value = 62760
$ws_BSM.Range("K84").Value = 90000
$ws_BSM.Range("L84").Value = 188280
$ws_BSM.Range("M84").Value = -84696
$ws_BSM.Range("N84").Value = -198888
$ws_BSM.Range("H99").Value = 2499
$ws_BSM.Range("I99").Value = 2499
$ws_BSM.Range("K99").Value = 2499
$ws_BSM.Range("M99").Value = -1001
$ws_BSM.Range("H105").Value = 3225.4
$ws_BSM.Range("I105").Value = 1454.6666
$ws_BSM.Range("K105").Value = 1454.6666
$ws_BSM.Range("M105").Value = 292.3334
$ws_BSM.Range("H134").Value = 28371.453
$ws_BSM.Range("I134").Value = 1441.8125
$ws_BSM.Range("K134").Value = 4325.4375
$ws_BSM.Range("M134").Value = -1790.4375

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 394026.47
$ws_CRP.Range("I31").Value = 1541458.9
$ws_CRP.Range("J31").Value = 11549
$ws_CRP.Range("K31").Value = 1541458.9
$ws_CRP.Range("L31").Value = 11549
$ws_CRP.Range("M31").Value = -1541163.9
$ws_CRP.Range("N31").Value = -12139
$ws_CRP.Range("H34").Value = 394026.47
$ws_CRP.Range("I34").Value = 1541458.9
$ws_CRP.Range("J34").Value = 11549
$ws_CRP.Range("K34").Value = 1541458.9
$ws_CRP.Range("L34").Value = 11549
$ws_CRP.Range("M34").Value = -1541256.9
$ws_CRP.Range("N34").Value = -11953
$ws_CRP.Range("H58").Value = 7449.6665
$ws_CRP.Range("I58").Value = 2171.2
$ws_CRP.Range("J58").Value = 12728.134
$ws_CRP.Range("K58").Value = 2171.2
$ws_CRP.Range("L58").Value = 12728.134
$ws_CRP.Range("M58").Value = -1968.2
$ws_CRP.Range("N58").Value = -13134.134
$ws_CRP.Range("H94").Value = 1059.8
$ws_CRP.Range("I94").Value = 899.5
$ws_CRP.Range("J94").Value = 1166.6666
$ws_CRP.Range("K94").Value = 899.5
$ws_CRP.Range("L94").Value = 1166.6666
$ws_CRP.Range("M94").Value = -448.5
$ws_CRP.Range("N94").Value = -2068.6666
$ws_CRP.Range("H107").Value = 1296.7428
$ws_CRP.Range("I107").Value = 989.8182
$ws_CRP.Range("J107").Value = 1816.1538
$ws_CRP.Range("K107").Value = 989.8182
$ws_CRP.Range("L107").Value = 1816.1538
$ws_CRP.Range("M107").Value = 930.1818
$ws_CRP.Range("N107").Value = -5656.1538
$ws_CRP.Range("H136").Value = 7449.6665
$ws_CRP.Range("I136").Value = 2171.2
$ws_CRP.Range("J136").Value = 12728.134
$ws_CRP.Range("K136").Value = 6513.599999999999
$ws_CRP.Range("L136").Value = 38184.402
$ws_CRP.Range("M136").Value = -3963.599999999999
$ws_CRP.Range("N136").Value = -43284.402

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 21027314
$ws_CUL.Range("I4").Value = 1233881.2
$ws_CUL.Range("K4").Value = 3701643.6
$ws_CUL.Range("M4").Value = -3701531.6
$ws_CUL.Range("H8").Value = 299.33334
$ws_CUL.Range("I8").Value = 299.33334
$ws_CUL.Range("K8").Value = 898.0000200000001
$ws_CUL.Range("M8").Value = -759.0000200000001
$ws_CUL.Range("H68").Value = 1951.6615
$ws_CUL.Range("I68").Value = 1666.5555
$ws_CUL.Range("J68").Value = 1997.4822
$ws_CUL.Range("K68").Value = 4999.666499999999
$ws_CUL.Range("L68").Value = 5992.446599999999
$ws_CUL.Range("M68").Value = -4188.666499999999
$ws_CUL.Range("N68").Value = -7614.446599999999
$ws_CUL.Range("H71").Value = 1951.6615
$ws_CUL.Range("I71").Value = 1666.5555
$ws_CUL.Range("J71").Value = 1997.4822
$ws_CUL.Range("K71").Value = 14998.9995
$ws_CUL.Range("L71").Value = 17977.3398
$ws_CUL.Range("M71").Value = -10942.9995
$ws_CUL.Range("N71").Value = -26089.3398
$ws_CUL.Range("H81").Value = 33601.867
$ws_CUL.Range("I81").Value = 1010.5
$ws_CUL.Range("J81").Value = 55329.445
$ws_CUL.Range("K81").Value = 3031.5
$ws_CUL.Range("L81").Value = 165988.335
$ws_CUL.Range("M81").Value = -1908.5
$ws_CUL.Range("N81").Value = -168234.335
$ws_CUL.Range("H84").Value = 33601.867
$ws_CUL.Range("I84").Value = 1010.5
$ws_CUL.Range("J84").Value = 55329.445
$ws_CUL.Range("K84").Value = 9094.5
$ws_CUL.Range("L84").Value = 497965.005
$ws_CUL.Range("M84").Value = -3478.5
$ws_CUL.Range("N84").Value = -509197.005
$ws_CUL.Range("H103").Value = 570.0909
$ws_CUL.Range("I103").Value = 216
$ws_CUL.Range("J103").Value = 995
$ws_CUL.Range("K103").Value = 648
$ws_CUL.Range("L103").Value = 2985
$ws_CUL.Range("M103").Value = 231
$ws_CUL.Range("N103").Value = -4743
$ws_CUL.Range("H106").Value = 10764.5
$ws_CUL.Range("I106").Value = 11500
$ws_CUL.Range("J106").Value = 10029
$ws_CUL.Range("K106").Value = 34500
$ws_CUL.Range("L106").Value = 30087
$ws_CUL.Range("M106").Value = -33554
$ws_CUL.Range("N106").Value = -31979
$ws_CUL.Range("H132").Value = 2861466.8
$ws_CUL.Range("J132").Value = 32052
$ws_CUL.Range("L132").Value = 288468
$ws_CUL.Range("N132").Value = -293528

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H97").Value = 1188.7222
$ws_GSM.Range("I97").Value = 1286.5333
$ws_GSM.Range("J97").Value = 699.6667
$ws_GSM.Range("K97").Value = 1286.5333
$ws_GSM.Range("L97").Value = 699.6667
$ws_GSM.Range("M97").Value = -790.5333000000001
$ws_GSM.Range("N97").Value = -1691.6667
$ws_GSM.Range("H132").Value = 86334.25
$ws_GSM.Range("I132").Value = 2574.625
$ws_GSM.Range("J132").Value = 253853.5
$ws_GSM.Range("K132").Value = 7723.875
$ws_GSM.Range("L132").Value = 761560.5
$ws_GSM.Range("M132").Value = -5193.875
$ws_GSM.Range("N132").Value = -766620.5
$ws_GSM.Range("H134").Value = 78773.336
$ws_GSM.Range("J134").Value = 78773.336
$ws_GSM.Range("L134").Value = 236320.008
$ws_GSM.Range("N134").Value = -241390.008

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H61").Value = 5668.7354
$ws_LTW.Range("I61").Value = 4822.857
$ws_LTW.Range("J61").Value = 6260.85
$ws_LTW.Range("K61").Value = 4822.857
$ws_LTW.Range("L61").Value = 6260.85
$ws_LTW.Range("M61").Value = -4620.857
$ws_LTW.Range("N61").Value = -6664.85
$ws_LTW.Range("H113").Value = 5668.7354
$ws_LTW.Range("I113").Value = 4822.857
$ws_LTW.Range("J113").Value = 6260.85
$ws_LTW.Range("K113").Value = 4822.857
$ws_LTW.Range("L113").Value = 6260.85
$ws_LTW.Range("M113").Value = -2652.857
$ws_LTW.Range("N113").Value = -10600.85
$ws_LTW.Range("H122").Value = 412428.25
$ws_LTW.Range("I122").Value = 3457.0667
$ws_LTW.Range("J122").Value = 923642.25
$ws_LTW.Range("K122").Value = 10371.2001
$ws_LTW.Range("L122").Value = 2770926.75
$ws_LTW.Range("M122").Value = -7921.2001
$ws_LTW.Range("N122").Value = -2775826.75
$ws_LTW.Range("H132").Value = 3365.5098
$ws_LTW.Range("I132").Value = 2753.1538
$ws_LTW.Range("K132").Value = 8259.4614
$ws_LTW.Range("M132").Value = -5729.4614
$ws_LTW.Range("H136").Value = 418764.75
$ws_LTW.Range("I136").Value = 562612.9399999999
$ws_LTW.Range("J136").Value = 11194.944
$ws_LTW.Range("K136").Value = 1687838.82
$ws_LTW.Range("L136").Value = 33584.83199999999
$ws_LTW.Range("M136").Value = -1685288.82
$ws_LTW.Range("N136").Value = -38684.83199999999

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H49").Value = 0
$ws_WVR.Range("I49").Value = 0
$ws_WVR.Range("K49").Value = 0
$ws_WVR.Range("H100").Value = 1495.7142
$ws_WVR.Range("I100").Value = 1661.6666
$ws_WVR.Range("K100").Value = 3323.3332
$ws_WVR.Range("M100").Value = -2782.3332
$ws_WVR.Range("H130").Value = 89995
$ws_WVR.Range("J130").Value = 89995
$ws_WVR.Range("L130").Value = 89995
$ws_WVR.Range("N130").Value = -100035
$ws_WVR.Range("H136").Value = 314174.6
$ws_WVR.Range("I136").Value = 336843.5
$ws_WVR.Range("K136").Value = 1010530.5
$ws_WVR.Range("M136").Value = -1007980.5

# M49 cell is removed entirely in the target (not just set to 0/blank)
$ws_WVR.Range("M49").ClearContents()
